# Scheduled market-data refresh: updates currentAveragePrice* / Leve*Price / Leve*Profit
# columns (H:N) for a handful of leve rows across several crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Cells.Item(6, 8).Value = 314.2143
$ws.Cells.Item(6, 9).Value = 139.9
$ws.Cells.Item(6, 10).Value = 750
$ws.Cells.Item(6, 11).Value = 419.7
$ws.Cells.Item(6, 12).Value = 2250
$ws.Cells.Item(6, 13).Value = -307.7
$ws.Cells.Item(6, 14).Value = -2474

# Row 116
$ws.Cells.Item(116, 8).Value = 6307.143
$ws.Cells.Item(116, 9).Value = 6854
$ws.Cells.Item(116, 11).Value = 6854
$ws.Cells.Item(116, 13).Value = -3412

# Row 132
$ws.Cells.Item(132, 8).Value = 7548.952
$ws.Cells.Item(132, 9).Value = 3008.5334
$ws.Cells.Item(132, 11).Value = 9025.600199999999
$ws.Cells.Item(132, 13).Value = -6495.600199999999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 517107.5
$ws.Cells.Item(61, 9).Value = 669100.8
$ws.Cells.Item(61, 10).Value = 422111.66
$ws.Cells.Item(61, 11).Value = 669100.8
$ws.Cells.Item(61, 12).Value = 422111.66
$ws.Cells.Item(61, 13).Value = -668888.8
$ws.Cells.Item(61, 14).Value = -422535.66

# Row 102
$ws.Cells.Item(102, 8).Value = 4343.7617
$ws.Cells.Item(102, 9).Value = 1401.1875
$ws.Cells.Item(102, 11).Value = 1401.1875
$ws.Cells.Item(102, 13).Value = 220.8125

# Row 136
$ws.Cells.Item(136, 8).Value = 517107.5
$ws.Cells.Item(136, 9).Value = 669100.8
$ws.Cells.Item(136, 10).Value = 422111.66
$ws.Cells.Item(136, 11).Value = 2007302.4
$ws.Cells.Item(136, 12).Value = 1266334.98
$ws.Cells.Item(136, 13).Value = -2004752.4
$ws.Cells.Item(136, 14).Value = -1271434.98

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1202
$ws.Cells.Item(20, 9).Value = 1027.6666
$ws.Cells.Item(20, 10).Value = 1439.7273
$ws.Cells.Item(20, 11).Value = 1027.6666
$ws.Cells.Item(20, 12).Value = 1439.7273
$ws.Cells.Item(20, 13).Value = -780.6666
$ws.Cells.Item(20, 14).Value = -1933.7273

# Row 107
$ws.Cells.Item(107, 8).Value = 1490
$ws.Cells.Item(107, 9).Value = 1290.3334
$ws.Cells.Item(107, 10).Value = 1789.5
$ws.Cells.Item(107, 11).Value = 1290.3334
$ws.Cells.Item(107, 12).Value = 1789.5
$ws.Cells.Item(107, 13).Value = 629.6666
$ws.Cells.Item(107, 14).Value = -5629.5

# Row 134
$ws.Cells.Item(134, 8).Value = 2003.3125
$ws.Cells.Item(134, 9).Value = 1105.3
$ws.Cells.Item(134, 10).Value = 3500
$ws.Cells.Item(134, 11).Value = 3315.9
$ws.Cells.Item(134, 12).Value = 10500
$ws.Cells.Item(134, 13).Value = -780.8999999999996
$ws.Cells.Item(134, 14).Value = -15570

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5534.5557
$ws.Cells.Item(31, 9).Value = 2204.5833
$ws.Cells.Item(31, 10).Value = 6745.4546
$ws.Cells.Item(31, 11).Value = 2204.5833
$ws.Cells.Item(31, 12).Value = 6745.4546
$ws.Cells.Item(31, 13).Value = -1909.5833
$ws.Cells.Item(31, 14).Value = -7335.4546

# Row 34
$ws.Cells.Item(34, 8).Value = 5534.5557
$ws.Cells.Item(34, 9).Value = 2204.5833
$ws.Cells.Item(34, 10).Value = 6745.4546
$ws.Cells.Item(34, 11).Value = 2204.5833
$ws.Cells.Item(34, 12).Value = 6745.4546
$ws.Cells.Item(34, 13).Value = -2002.5833
$ws.Cells.Item(34, 14).Value = -7149.4546

# Row 86
$ws.Cells.Item(86, 8).Value = 3571.55
$ws.Cells.Item(86, 9).Value = 2617.5908
$ws.Cells.Item(86, 10).Value = 4737.5
$ws.Cells.Item(86, 11).Value = 2617.5908
$ws.Cells.Item(86, 12).Value = 4737.5
$ws.Cells.Item(86, 13).Value = -1494.5908
$ws.Cells.Item(86, 14).Value = -6983.5

# Row 89
$ws.Cells.Item(89, 8).Value = 3571.55
$ws.Cells.Item(89, 9).Value = 2617.5908
$ws.Cells.Item(89, 10).Value = 4737.5
$ws.Cells.Item(89, 11).Value = 13087.954
$ws.Cells.Item(89, 12).Value = 23687.5
$ws.Cells.Item(89, 13).Value = -7471.954
$ws.Cells.Item(89, 14).Value = -34919.5

# Row 134
$ws.Cells.Item(134, 8).Value = 1696.6364
$ws.Cells.Item(134, 9).Value = 1125.0625
$ws.Cells.Item(134, 10).Value = 2234.5881
$ws.Cells.Item(134, 11).Value = 3375.1875
$ws.Cells.Item(134, 12).Value = 6703.7643
$ws.Cells.Item(134, 13).Value = -840.1875
$ws.Cells.Item(134, 14).Value = -11773.7643

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1360.4231
$ws.Cells.Item(5, 9).Value = 828.2857
$ws.Cells.Item(5, 10).Value = 1981.25
$ws.Cells.Item(5, 11).Value = 2484.8571
$ws.Cells.Item(5, 12).Value = 5943.75
$ws.Cells.Item(5, 13).Value = -2372.8571
$ws.Cells.Item(5, 14).Value = -6167.75

# Row 68
$ws.Cells.Item(68, 8).Value = 558.3333
$ws.Cells.Item(68, 10).Value = 650
$ws.Cells.Item(68, 12).Value = 1950
$ws.Cells.Item(68, 14).Value = -3572

# Row 71
$ws.Cells.Item(71, 8).Value = 558.3333
$ws.Cells.Item(71, 10).Value = 650
$ws.Cells.Item(71, 12).Value = 5850
$ws.Cells.Item(71, 14).Value = -13962

# Row 86
$ws.Cells.Item(86, 8).Value = 626.5333000000001
$ws.Cells.Item(86, 9).Value = 641.5833
$ws.Cells.Item(86, 10).Value = 566.3333
$ws.Cells.Item(86, 11).Value = 1924.7499
$ws.Cells.Item(86, 12).Value = 1698.9999
$ws.Cells.Item(86, 13).Value = -738.7499
$ws.Cells.Item(86, 14).Value = -4070.9999

# Row 89
$ws.Cells.Item(89, 8).Value = 626.5333000000001
$ws.Cells.Item(89, 9).Value = 641.5833
$ws.Cells.Item(89, 10).Value = 566.3333
$ws.Cells.Item(89, 11).Value = 5774.2497
$ws.Cells.Item(89, 12).Value = 5096.9997
$ws.Cells.Item(89, 13).Value = 153.7502999999997
$ws.Cells.Item(89, 14).Value = -16952.9997

# Row 92
$ws.Cells.Item(92, 8).Value = 600
$ws.Cells.Item(92, 9).Value = 450
$ws.Cells.Item(92, 11).Value = 1350
$ws.Cells.Item(92, 13).Value = -102

# Row 122
$ws.Cells.Item(122, 8).Value = 1032.4546
$ws.Cells.Item(122, 9).Value = 359.14285
$ws.Cells.Item(122, 10).Value = 1213.7307
$ws.Cells.Item(122, 11).Value = 3232.28565
$ws.Cells.Item(122, 12).Value = 10923.5763
$ws.Cells.Item(122, 13).Value = -782.2856500000003
$ws.Cells.Item(122, 14).Value = -15823.5763

# Row 135
$ws.Cells.Item(135, 8).Value = 1360.4231
$ws.Cells.Item(135, 9).Value = 828.2857
$ws.Cells.Item(135, 10).Value = 1981.25
$ws.Cells.Item(135, 11).Value = 7454.571300000001
$ws.Cells.Item(135, 12).Value = 17831.25
$ws.Cells.Item(135, 13).Value = -4919.571300000001
$ws.Cells.Item(135, 14).Value = -22901.25

# Row 139
$ws.Cells.Item(139, 8).Value = 2494.375
$ws.Cells.Item(139, 9).Value = 1109.1666
$ws.Cells.Item(139, 10).Value = 6650
$ws.Cells.Item(139, 11).Value = 3327.4998
$ws.Cells.Item(139, 12).Value = 19950
$ws.Cells.Item(139, 13).Value = 1812.5002
$ws.Cells.Item(139, 14).Value = -30230

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 3315.44
$ws.Cells.Item(102, 9).Value = 1493.3529
$ws.Cells.Item(102, 10).Value = 7187.375
$ws.Cells.Item(102, 11).Value = 1493.3529
$ws.Cells.Item(102, 12).Value = 7187.375
$ws.Cells.Item(102, 13).Value = 128.6470999999999
$ws.Cells.Item(102, 14).Value = -10431.375

# Row 107
$ws.Cells.Item(107, 8).Value = 6207.8823
$ws.Cells.Item(107, 10).Value = 546
$ws.Cells.Item(107, 12).Value = 546
$ws.Cells.Item(107, 14).Value = -4386

# Row 126
$ws.Cells.Item(126, 8).Value = 2715.0454
$ws.Cells.Item(126, 9).Value = 2001.1
$ws.Cells.Item(126, 10).Value = 3310
$ws.Cells.Item(126, 11).Value = 6003.299999999999
$ws.Cells.Item(126, 12).Value = 9930
$ws.Cells.Item(126, 13).Value = -3533.299999999999
$ws.Cells.Item(126, 14).Value = -14870

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 195.64706
$ws.Cells.Item(55, 9).Value = 166.09091
$ws.Cells.Item(55, 10).Value = 249.83333
$ws.Cells.Item(55, 11).Value = 166.09091
$ws.Cells.Item(55, 12).Value = 249.83333
$ws.Cells.Item(55, 13).Value = 6.909089999999992
$ws.Cells.Item(55, 14).Value = -595.8333299999999

# Row 122
$ws.Cells.Item(122, 8).Value = 1711.45
$ws.Cells.Item(122, 9).Value = 1687.2307
$ws.Cells.Item(122, 10).Value = 1756.4286
$ws.Cells.Item(122, 11).Value = 5061.6921
$ws.Cells.Item(122, 12).Value = 5269.2858
$ws.Cells.Item(122, 13).Value = -2611.6921
$ws.Cells.Item(122, 14).Value = -10169.2858

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 2974
$ws.Cells.Item(132, 9).Value = 2352.2354
$ws.Cells.Item(132, 11).Value = 7056.706200000001
$ws.Cells.Item(132, 13).Value = -4526.706200000001

# Row 136
$ws.Cells.Item(136, 8).Value = 27326988
$ws.Cells.Item(136, 9).Value = 41710690
$ws.Cells.Item(136, 10).Value = 772461.9399999999
$ws.Cells.Item(136, 11).Value = 125132070
$ws.Cells.Item(136, 12).Value = 2317385.82
$ws.Cells.Item(136, 13).Value = -125129520
$ws.Cells.Item(136, 14).Value = -2322485.82
